$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.925.86'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '2.222.01'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '291.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.21%  '

$ws.Range("E13").Value = '  +2.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").Value = '2.563.24'
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.40%  '

$ws.Range("D17").Value = '2.218.67'
$ws.Range("E17").Value = '  -1.14%  '

$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("D19").Value = '39.876.68'
$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.29%  '

$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.64%  '

$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("E30").Value = '  -6.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.28%  '

$ws.Range("E33").Value = '  +0.13%  '

$ws.Range("E34").Value = '  +0.59%  '

$ws.Range("E35").Value = '  +7.19%  '

$ws.Range("E36").Value = '  -0.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.112'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0991'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.10%  '

$ws.Range("D42").Value = '2.107.24'
$ws.Range("E42").Value = '  +2.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.36%  '

$ws.Range("E45").Value = '  +0.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("E47").Value = '  -4.85%  '

$ws.Range("E48").Value = '  +3.64%  '

$ws.Range("D49").Value = '2.434.88'
$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E50").Value = '  +2.89%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.40%  '
